$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column order for the header row (row 2), per the "re-ordered columns" commit.
# IndicatorValues_ID and the *BusinessKey columns move to the front (alphabetically
# sorted among themselves, with IndicatorValues_ID first), followed by the
# remaining original columns in their prior relative order.
$newOrder = @(
    "IndicatorValues_ID",
    "AgeBandBusinessKey",
    "BusinessKey",
    "CommunityTypeBusinessKey",
    "DataVersionBusinessKey",
    "DonorBusinessKey",
    "FrameworkBusinessKey",
    "GenderBusinessKey",
    "GroupBusinessKey",
    "IndicatorBusinessKey",
    "InstitutionBusinessKey",
    "LocationBusinessKey",
    "OrganizationBusinessKey",
    "ReportingPeriodBusinessKey",
    "ResultAreaBusinessKey",
    "StrategicElementBusinessKey",
    "ActualDate",
    "ActualLabel",
    "ActualValue",
    "GroupVersion",
    "Notes"
)

for ($i = 0; $i -lt $newOrder.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $newOrder[$i]
}
